$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 94
$rng = $ws.Range("A$row`:J$row")

# Force these cells to store the values as literal text (matching the
# source file, which keeps every value - even numeric-looking ones - as
# an inline/shared string) instead of letting Excel auto-convert them to
# numbers or dates.
$rng.NumberFormat = "@"

$ws.Range("A$row").Value = "2025-06-03"
$ws.Range("B$row").Value = "35.5"
$ws.Range("C$row").Value = "35.01"
$ws.Range("D$row").Value = "0.94"
$ws.Range("E$row").Value = "0.253"
$ws.Range("F$row").Value = "0.09"
$ws.Range("G$row").Value = "5,352"
$ws.Range("H$row").Value = "8,013"
$ws.Range("I$row").Value = "8,063"
$ws.Range("J$row").Value = "7.2065"

# Revert the cell style back to the workbook's default ("Normal") so no
# new cell-format style gets introduced - the text values already typed
# above stay as text even once the display format goes back to General.
$rng.Style = "Normal"
